# Update "want to go" counts (column F) for a handful of events on the
# "展览" (Exhibitions) and "全部类型" (All types) sheets, reflecting refreshed
# scrape output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 11279
$wsExhibit.Range("F9").Value  = 11200
$wsExhibit.Range("F10").Value = 455
$wsExhibit.Range("F11").Value = 1146
$wsExhibit.Range("F14").Value = 5603
$wsExhibit.Range("F15").Value = 99
$wsExhibit.Range("F16").Value = 3457

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 11280
$wsAll.Range("F11").Value = 11200
$wsAll.Range("F12").Value = 455
$wsAll.Range("F13").Value = 1146
$wsAll.Range("F16").Value = 5603
$wsAll.Range("F17").Value = 99
$wsAll.Range("F18").Value = 3457
